# Apply performance-fixture data refresh to both worksheets.
$wb = $excel.ActiveWorkbook

$wsWarm = $wb.Worksheets.Item("Warmlaufen")
$wsPerf = $wb.Worksheets.Item("Performanz Messung")

# Update the timestamp title shown in cell A1 on both sheets (shared text).
$wsWarm.Range("A1").Value = "05.05.2016 um 21:52 Uhr"
$wsPerf.Range("A1").Value = "05.05.2016 um 21:52 Uhr"

# --- Sheet "Warmlaufen" ---
$wsWarm.Range("C3").Value = 3.0
$wsWarm.Range("F3").Value = 0.3
$wsWarm.Range("G3").Value = 0.674949

$wsWarm.Range("C4").Value = 630.0
$wsWarm.Range("D4").Value = 26.0
$wsWarm.Range("E4").Value = 212.0
$wsWarm.Range("F4").Value = 63.0
$wsWarm.Range("G4").Value = 56.1169

$wsWarm.Range("C5").Value = 129.0
$wsWarm.Range("D5").Value = 1.0
$wsWarm.Range("E5").Value = 88.0
$wsWarm.Range("F5").Value = 12.9
$wsWarm.Range("G5").Value = 26.451
$wsWarm.Range("H5").Value = 1.0

$wsWarm.Range("C6").Value = 13.0
$wsWarm.Range("E6").Value = 3.0
$wsWarm.Range("F6").Value = 1.3
$wsWarm.Range("G6").Value = 0.823273

$wsWarm.Range("C7").Value = 508.0
$wsWarm.Range("D7").Value = 15.0
$wsWarm.Range("E7").Value = 256.0
$wsWarm.Range("F7").Value = 50.8
$wsWarm.Range("G7").Value = 73.3285
$wsWarm.Range("H7").Value = 23.0

$wsWarm.Range("D8").Value = 0.0
$wsWarm.Range("E8").Value = 2.0
$wsWarm.Range("G8").Value = 0.471405
$wsWarm.Range("H8").Value = 0.0

# --- Sheet "Performanz Messung" ---
$wsPerf.Range("B3").Value = 32939.0
$wsPerf.Range("C3").Value = 541.0
$wsPerf.Range("E3").Value = 6.0
$wsPerf.Range("F3").Value = 0.0164243
$wsPerf.Range("G3").Value = 0.133396

$wsPerf.Range("B4").Value = 32939.0
$wsPerf.Range("C4").Value = 366885.0
$wsPerf.Range("E4").Value = 57.0
$wsPerf.Range("F4").Value = 11.1383
$wsPerf.Range("G4").Value = 1.32497

$wsPerf.Range("B5").Value = 32939.0
$wsPerf.Range("C5").Value = 17482.0
$wsPerf.Range("E5").Value = 37.0
$wsPerf.Range("F5").Value = 0.530739
$wsPerf.Range("G5").Value = 0.557399
$wsPerf.Range("H5").Value = 1.0

$wsPerf.Range("B6").Value = 32939.0
$wsPerf.Range("C6").Value = 5700.0
$wsPerf.Range("E6").Value = 12.0
$wsPerf.Range("F6").Value = 0.173047
$wsPerf.Range("G6").Value = 0.388118

$wsPerf.Range("B7").Value = 32939.0
$wsPerf.Range("C7").Value = 185247.0
$wsPerf.Range("E7").Value = 31.0
$wsPerf.Range("F7").Value = 5.62394
$wsPerf.Range("G7").Value = 0.964435

$wsPerf.Range("B8").Value = 32939.0
$wsPerf.Range("C8").Value = 493.0
$wsPerf.Range("F8").Value = 0.0149671
$wsPerf.Range("G8").Value = 0.122667

$wsPerf.Range("C12").Value = 0.03
$wsPerf.Range("D12").Value = 0.016
$wsPerf.Range("E12").Value = 0.0139
$wsPerf.Range("F12").Value = 0.01595
$wsPerf.Range("G12").Value = 0.0159

$wsPerf.Range("B13").Value = 26.9
$wsPerf.Range("C13").Value = 16.85
$wsPerf.Range("D13").Value = 11.912
$wsPerf.Range("E13").Value = 11.2144
$wsPerf.Range("F13").Value = 11.0891
$wsPerf.Range("G13").Value = 11.1288

$wsPerf.Range("B14").Value = 4.9
$wsPerf.Range("C14").Value = 1.16
$wsPerf.Range("D14").Value = 0.633
$wsPerf.Range("E14").Value = 0.5478
$wsPerf.Range("F14").Value = 0.53215
$wsPerf.Range("G14").Value = 0.531

$wsPerf.Range("B15").Value = 0.7
$wsPerf.Range("C15").Value = 0.33
$wsPerf.Range("D15").Value = 0.249
$wsPerf.Range("E15").Value = 0.1869
$wsPerf.Range("F15").Value = 0.1731
$wsPerf.Range("G15").Value = 0.172167

$wsPerf.Range("B16").Value = 14.1
$wsPerf.Range("D16").Value = 6.108
$wsPerf.Range("E16").Value = 5.6631
$wsPerf.Range("F16").Value = 5.5901
$wsPerf.Range("G16").Value = 5.61677

$wsPerf.Range("B17").Value = 0.8
$wsPerf.Range("C17").Value = 0.61
$wsPerf.Range("D17").Value = 0.091
$wsPerf.Range("E17").Value = 0.0216
$wsPerf.Range("F17").Value = 0.01685
$wsPerf.Range("G17").Value = 0.0153

$wb.Save()
